$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B3").Value = 57073
$ws.Range("B4").Value = 57073
$ws.Range("B5").Value = 57073
$ws.Range("B6").Value = 57073
